$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "newModel"
$ws.Range("D1").Value = "newModel2"
$ws.Range("E1").Value = "newModel3"
$ws.Range("F1").Value = "newModel7"
$ws.Range("G1").Value = "newModel1"

$ws.Range("C2:G2").Value = 1
$ws.Range("C6:G6").Value = 1
$ws.Range("C8:G8").Value = 1
$ws.Range("C10:G10").Value = 1
